# Apply the "Fixed update to excel issue" edit:
#  1. Rename the "Requested quantity" header on "Weekly Quantity" sheet to "Weekly_PO_Qty"
#  2. Rename the "Requested quantity" header on "Monthly Trend" sheet to "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Rename header in "Weekly Quantity" sheet ---
$ws1.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Rename header in "Monthly Trend" sheet ---
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Create the new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Replicate the header (bold/centered) and date-number-format styling that is
# used on the other two sheets, so the new sheet matches the workbook look and feel.
$ws1.Range("A1").Copy() | Out-Null
$ws3.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Copy() | Out-Null
$ws3.Range("A2:A15").PasteSpecial(-4122) | Out-Null

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Data rows
$ws3.Range("A2").Value = 45298.99999999999
$ws3.Range("B2").Value = 23
$ws3.Range("C2").Value = -14.49030783587633
$ws3.Range("D2").Value = 58.76085296754226

$ws3.Range("A3").Value = 45333.99999999999
$ws3.Range("B3").Value = 28
$ws3.Range("C3").Value = -8.744382106312658
$ws3.Range("D3").Value = 65.35801653637775

$ws3.Range("A4").Value = 45438.99999999999
$ws3.Range("B4").Value = 44
$ws3.Range("C4").Value = 7.031220819028939
$ws3.Range("D4").Value = 80.6169372037256

$ws3.Range("A5").Value = 45445.99999999999
$ws3.Range("B5").Value = 45
$ws3.Range("C5").Value = 9.074445377945349
$ws3.Range("D5").Value = 83.49402246792513

$ws3.Range("A6").Value = 45452.99999999999
$ws3.Range("B6").Value = 46
$ws3.Range("C6").Value = 6.416116852768758
$ws3.Range("D6").Value = 80.52435112132835

$ws3.Range("A7").Value = 45459.99999999999
$ws3.Range("B7").Value = 47
$ws3.Range("C7").Value = 11.16276316575149
$ws3.Range("D7").Value = 83.71450895482295

$ws3.Range("A8").Value = 45466.99999999999
$ws3.Range("B8").Value = 48
$ws3.Range("C8").Value = 13.44969917201085
$ws3.Range("D8").Value = 84.90181384875838

$ws3.Range("A9").Value = 45473.99999999999
$ws3.Range("B9").Value = 49
$ws3.Range("C9").Value = 12.34037826524282
$ws3.Range("D9").Value = 86.64763177062878

$ws3.Range("A10").Value = 45480.99999999999
$ws3.Range("B10").Value = 50
$ws3.Range("C10").Value = 13.25092080184036
$ws3.Range("D10").Value = 87.11811898859216

$ws3.Range("A11").Value = 45487.99999999999
$ws3.Range("B11").Value = 51
$ws3.Range("C11").Value = 13.72709386422437
$ws3.Range("D11").Value = 87.82846363111928

$ws3.Range("A12").Value = 45494.99999999999
$ws3.Range("B12").Value = 52
$ws3.Range("C12").Value = 14.38126913177011
$ws3.Range("D12").Value = 89.7905331783107

$ws3.Range("A13").Value = 45501.99999999999
$ws3.Range("B13").Value = 53
$ws3.Range("C13").Value = 14.89949512499028
$ws3.Range("D13").Value = 92.14689084923822

$ws3.Range("A14").Value = 45508.99999999999
$ws3.Range("B14").Value = 54
$ws3.Range("C14").Value = 18.32990386292177
$ws3.Range("D14").Value = 89.13157014150995

$ws3.Range("A15").Value = 45515.99999999999
$ws3.Range("B15").Value = 55
$ws3.Range("C15").Value = 15.34747428126287
$ws3.Range("D15").Value = 90.60642675394641

# Leave the active cell/selection at A1 on the new sheet, matching the source workbook convention
$ws3.Range("A1").Select() | Out-Null

# Restore the originally active sheet (first sheet) so the workbook-level view state is unchanged
$ws1.Select() | Out-Null
